# Auto-generated script to apply numeric updates described in the commit diff.
# The workbook contains plain numeric leveling-profit data (no formulas).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 8703.291999999999
$ws.Range("I76").Value = 13747.9
$ws.Range("K76").Value = 13747.9
$ws.Range("M76").Value = -13432.9
$ws.Range("H79").Value = 8703.291999999999
$ws.Range("I79").Value = 13747.9
$ws.Range("K79").Value = 13747.9
$ws.Range("M79").Value = -12655.9
$ws.Range("H86").Value = 4088.6553
$ws.Range("I86").Value = 3054
$ws.Range("J86").Value = 4554.25
$ws.Range("K86").Value = 3054
$ws.Range("L86").Value = 4554.25
$ws.Range("M86").Value = -1931
$ws.Range("N86").Value = -6800.25
$ws.Range("H88").Value = 4297.95
$ws.Range("I88").Value = 521.4286
$ws.Range("K88").Value = 521.4286
$ws.Range("M88").Value = -115.4286
$ws.Range("H89").Value = 4088.6553
$ws.Range("I89").Value = 3054
$ws.Range("J89").Value = 4554.25
$ws.Range("K89").Value = 15270
$ws.Range("L89").Value = 22771.25
$ws.Range("M89").Value = -9654
$ws.Range("N89").Value = -34003.25
$ws.Range("H91").Value = 4297.95
$ws.Range("I91").Value = 521.4286
$ws.Range("K91").Value = 521.4286
$ws.Range("M91").Value = 882.5714
$ws.Range("H113").Value = 3714.1
$ws.Range("I113").Value = 3109
$ws.Range("J113").Value = 4176.8237
$ws.Range("K113").Value = 3109
$ws.Range("L113").Value = 4176.8237
$ws.Range("M113").Value = 145
$ws.Range("N113").Value = -10684.8237
$ws.Range("H116").Value = 86951.48
$ws.Range("I116").Value = 107319.35
$ws.Range("K116").Value = 107319.35
$ws.Range("M116").Value = -103877.35
$ws.Range("H132").Value = 3141.717
$ws.Range("I132").Value = 1519.4043
$ws.Range("K132").Value = 4558.2129
$ws.Range("M132").Value = -2028.2129
$ws.Range("H135").Value = 573
$ws.Range("I135").Value = 492.6
$ws.Range("K135").Value = 4433.400000000001
$ws.Range("M135").Value = -1898.400000000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4727.25
$ws.Range("I32").Value = 3271.5474
$ws.Range("J32").Value = 32385.6
$ws.Range("K32").Value = 3271.5474
$ws.Range("L32").Value = 32385.6
$ws.Range("M32").Value = -2984.5474
$ws.Range("N32").Value = -32959.6
$ws.Range("H38").Value = 1954.75
$ws.Range("I38").Value = 1954.75
$ws.Range("K38").Value = 1954.75
$ws.Range("M38").Value = -1487.75
$ws.Range("H61").Value = 2599.1865
$ws.Range("I61").Value = 1727.5
$ws.Range("J61").Value = 5156.1333
$ws.Range("K61").Value = 1727.5
$ws.Range("L61").Value = 5156.1333
$ws.Range("M61").Value = -1515.5
$ws.Range("N61").Value = -5580.1333
$ws.Range("H88").Value = 3000
$ws.Range("J88").Value = 2044.4445
$ws.Range("L88").Value = 2044.4445
$ws.Range("N88").Value = -2856.4445
$ws.Range("H91").Value = 3000
$ws.Range("J91").Value = 2044.4445
$ws.Range("L91").Value = 2044.4445
$ws.Range("N91").Value = -4852.4445
$ws.Range("H97").Value = 2450.8333
$ws.Range("J97").Value = 1766.6666
$ws.Range("L97").Value = 1766.6666
$ws.Range("N97").Value = -2758.6666
$ws.Range("H114").Value = 34333.332
$ws.Range("J114").Value = 34333.332
$ws.Range("L114").Value = 34333.332
$ws.Range("N114").Value = -43011.332
$ws.Range("H124").Value = 12612.7
$ws.Range("J124").Value = 12612.7
$ws.Range("L124").Value = 12612.7
$ws.Range("N124").Value = -22432.7
$ws.Range("H136").Value = 2599.1865
$ws.Range("I136").Value = 1727.5
$ws.Range("J136").Value = 5156.1333
$ws.Range("K136").Value = 5182.5
$ws.Range("L136").Value = 15468.3999
$ws.Range("M136").Value = -2632.5
$ws.Range("N136").Value = -20568.3999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 512.5
$ws.Range("I15").Value = 150
$ws.Range("J15").Value = 545.4545000000001
$ws.Range("K15").Value = 150
$ws.Range("L15").Value = 545.4545000000001
$ws.Range("M15").Value = 20
$ws.Range("N15").Value = -885.4545000000001
$ws.Range("H38").Value = 7032.2856
$ws.Range("I38").Value = 2900
$ws.Range("J38").Value = 7721
$ws.Range("K38").Value = 2900
$ws.Range("L38").Value = 7721
$ws.Range("M38").Value = -2523
$ws.Range("N38").Value = -8475
$ws.Range("H46").Value = 7032.2856
$ws.Range("I46").Value = 2900
$ws.Range("J46").Value = 7721
$ws.Range("K46").Value = 2900
$ws.Range("L46").Value = 7721
$ws.Range("M46").Value = -2689
$ws.Range("N46").Value = -8143
$ws.Range("H107").Value = 302.96
$ws.Range("I107").Value = 238.29411
$ws.Range("J107").Value = 440.375
$ws.Range("K107").Value = 238.29411
$ws.Range("L107").Value = 440.375
$ws.Range("M107").Value = 1681.70589
$ws.Range("N107").Value = -4280.375

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 53
$ws.Range("I10").Value = 53
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 159
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -20
$ws.Range("N10").ClearContents()
$ws.Range("H122").Value = 560.95
$ws.Range("I122").Value = 377.70587
$ws.Range("K122").Value = 3399.35283
$ws.Range("M122").Value = -949.3528299999998
$ws.Range("H132").Value = 3814.9487
$ws.Range("I132").Value = 2171.0625
$ws.Range("J132").Value = 4958.522
$ws.Range("K132").Value = 19539.5625
$ws.Range("L132").Value = 44626.698
$ws.Range("M132").Value = -17009.5625
$ws.Range("N132").Value = -49686.698

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 8475
$ws.Range("I10").Value = 8475
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 8475
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -8306
$ws.Range("N10").ClearContents()
$ws.Range("H123").Value = 19683.334
$ws.Range("J123").Value = 19683.334
$ws.Range("L123").Value = 19683.334
$ws.Range("N123").Value = -24583.334
$ws.Range("H126").Value = 3161.1
$ws.Range("I126").Value = 3451.375
$ws.Range("K126").Value = 10354.125
$ws.Range("M126").Value = -7884.125

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4422.222
$ws.Range("I136").Value = 2362.8147
$ws.Range("J136").Value = 7511.3335
$ws.Range("K136").Value = 7088.4441
$ws.Range("L136").Value = 22534.0005
$ws.Range("M136").Value = -4538.4441
$ws.Range("N136").Value = -27634.0005

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 300
$ws.Range("I7").Value = 300
$ws.Range("K7").Value = 300
$ws.Range("M7").Value = -187
$ws.Range("H107").Value = 880.2727
$ws.Range("I107").Value = 1500
$ws.Range("J107").Value = 647.875
$ws.Range("K107").Value = 4500
$ws.Range("L107").Value = 1943.625
$ws.Range("M107").Value = -2580
$ws.Range("N107").Value = -5783.625
$ws.Range("H122").Value = 85659.25
$ws.Range("I122").Value = 102091.1
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 306273.3
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -303823.3
$ws.Range("N122").Value = -15400
$ws.Range("H136").Value = 27029458
$ws.Range("I136").Value = 71431030
$ws.Range("J136").Value = 2413.4783
$ws.Range("K136").Value = 214293090
$ws.Range("L136").Value = 7240.4349
$ws.Range("M136").Value = -214290540
$ws.Range("N136").Value = -12340.4349
